# Update "想去人数" (number of people interested) figures on the
# "展览" (Exhibition) and "全部类型" (All types) sheets to match the
# newly scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 111
$ws1.Range("F5").Value = 5024
$ws1.Range("F6").Value = 374
$ws1.Range("F7").Value = 626
$ws1.Range("F8").Value = 294
$ws1.Range("F9").Value = 759
$ws1.Range("F10").Value = 243
$ws1.Range("F11").Value = 5

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 111
$ws4.Range("F5").Value = 5024
$ws4.Range("F6").Value = 374
$ws4.Range("F7").Value = 626
$ws4.Range("F8").Value = 294
$ws4.Range("F9").Value = 759
$ws4.Range("F10").Value = 30
$ws4.Range("F11").Value = 243
$ws4.Range("F12").Value = 5
